$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '19.916.73'
Set-TextValue 'E2' '  -5.90%  '
Set-TextValue 'D3' '1.408.82'
Set-TextValue 'E3' '  -6.99%  '
Set-TextValue 'D4' '1.003'
Set-TextValue 'E4' '  -0.46%  '
Set-TextValue 'E5' '  -0.46%  '
Set-TextValue 'D6' '275.80'
Set-TextValue 'E6' '  -3.70%  '
Set-TextValue 'D7' '0.3660'
Set-TextValue 'E7' '  -5.97%  '
Set-TextValue 'D8' '0.3094'
Set-TextValue 'E8' '  -1.47%  '
Set-TextValue 'D9' '39.74'
Set-TextValue 'E9' '  -6.00%  '
Set-TextValue 'D10' '1.032'
Set-TextValue 'E10' '  -1.67%  '
Set-TextValue 'D11' '0.06506'
Set-TextValue 'E11' '  -7.38%  '
Set-TextValue 'E12' '  -0.49%  '
Set-TextValue 'D13' '5.472'
Set-TextValue 'E13' '  -3.09%  '
Set-TextValue 'D14' '17.53'
Set-TextValue 'E14' '  -2.25%  '
Set-TextValue 'D15' '6.174'
Set-TextValue 'E15' '  -3.26%  '
Set-TextValue 'D16' '1.410.53'
Set-TextValue 'E16' '  -7.18%  '
Set-TextValue 'D17' '0.00001016'
Set-TextValue 'E17' '  -5.59%  '
Set-TextValue 'D18' '0.05664'
Set-TextValue 'E18' '  -13.84%  '
Set-TextValue 'E19' '  -0.44%  '
Set-TextValue 'D20' '70.77'
Set-TextValue 'E20' '  -13.92%  '
Set-TextValue 'D21' '5.609'
Set-TextValue 'E21' '  -6.74%  '
Set-TextValue 'D22' '14.68'
Set-TextValue 'E22' '  -3.35%  '
Set-TextValue 'D23' '10.87'
Set-TextValue 'E23' '  +0.94%  '
Set-TextValue 'D24' '2.236'
Set-TextValue 'E24' '  -4.85%  '
Set-TextValue 'D25' '19.924.54'
Set-TextValue 'E25' '  -5.86%  '
Set-TextValue 'D26' '2.251'
Set-TextValue 'E26' '  -4.69%  '
Set-TextValue 'D27' '132.71'
Set-TextValue 'E27' '  -10.18%  '
Set-TextValue 'E28' '  -4.65%  '
Set-TextValue 'D29' '1.569.14'
Set-TextValue 'E29' '  -7.14%  '
Set-TextValue 'D30' '109.64'
Set-TextValue 'E30' '  -4.80%  '
Set-TextValue 'D31' '3.913'
Set-TextValue 'E31' '  -18.53%  '
Set-TextValue 'D32' '5.251'
Set-TextValue 'E32' '  -12.12%  '
Set-TextValue 'D33' '0.8137'
Set-TextValue 'E33' '  -14.69%  '
Set-TextValue 'D34' '0.07670'
Set-TextValue 'E34' '  -4.02%  '
Set-TextValue 'D35' '1.483'
Set-TextValue 'E35' '  +0.11%  '
Set-TextValue 'D36' '8.290'
Set-TextValue 'E36' '  -1.69%  '
Set-TextValue 'D37' '4.896'
Set-TextValue 'E37' '  -3.66%  '
Set-TextValue 'D38' '0.05768'
Set-TextValue 'E38' '  -1.23%  '
Set-TextValue 'E39' '  -0.57%  '
Set-TextValue 'D40' '0.02055'
Set-TextValue 'E40' '  -4.28%  '
Set-TextValue 'D41' '10.42'
Set-TextValue 'E41' '  -7.80%  '
Set-TextValue 'D42' '0.1886'
Set-TextValue 'E42' '  -5.15%  '
Set-TextValue 'D43' '1.090'
Set-TextValue 'E43' '  -6.07%  '
Set-TextValue 'D44' '0.5296'
Set-TextValue 'E44' '  -6.27%  '
Set-TextValue 'B45' 'PancakeSwap'
Set-TextValue 'C45' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D45' '3.536'
Set-TextValue 'E45' '  -4.42%  '
Set-TextValue 'B46' 'EnergySwap'
Set-TextValue 'C46' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D46' '12.27'
Set-TextValue 'E46' '  -5.37%  '
Set-TextValue 'D47' '0.5172'
Set-TextValue 'E47' '  -5.45%  '
Set-TextValue 'D48' '114.21'
Set-TextValue 'E48' '  -0.36%  '
Set-TextValue 'D49' '1.763'
Set-TextValue 'E49' '  -4.84%  '
Set-TextValue 'D50' '1.031'
Set-TextValue 'E50' '  -9.72%  '
Set-TextValue 'D51' '1.003'
Set-TextValue 'E51' '  -0.45%  '
